$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right 140 -> 112, Max label "140 / 140" -> "112 / 112"
$ws.Range("B12").Value = 112
$ws.Range("E12").Value = "112 / 112"
